$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 14-25: update C (ID_db_Titolo_documento_fonte) and E (Titolo_documento_fonte) ---

$ws.Range("C14").Value = 8000035
$ws.Range("E14").Value = "Elemento d'intervento e di supporto dello Stato maggiore federale Protezione della popolazione"

$ws.Range("C15").Value = 8000036
$ws.Range("E15").Value = "Legge sulla protezione civile del 26 febbraio 2007"

$ws.Range("E16").Value = "Il comando della protezione civile"

$ws.Range("C17").Value = 8000039
$ws.Range("E17").Value = "Costruzioni di protezione"

$ws.Range("C18").Value = 8000042
$ws.Range("E18").Value = "Legge federale sulla protezione della popolazione e sulla protezione civile del 4 ottobre 2002"

$ws.Range("C19").Value = 8000047
$ws.Range("E19").Value = "Legge sulla protezione della popolazione (del 26 febbraio 2007)"

$ws.Range("C20").Value = 8000054
$ws.Range("E20").Value = "Regolamento sulla protezione della popolazione (RProtPop) (del 18 ottobre 2017)"

$ws.Range("C21").Value = 8000058
$ws.Range("E21").Value = "Segnali di allarme in Svizzera"

$ws.Range("C22").Value = 8000063
$ws.Range("E22").Value = "Dipartimento"

$ws.Range("C23").Value = 8000073
$ws.Range("E23").Value = "Le SOREU"

$ws.Range("C24").Value = 8000074
$ws.Range("E24").Value = "SOREU dei Laghi"

$ws.Range("C25").Value = 8000075
$ws.Range("E25").Value = "Chi siamo"

# --- Rows 32-40: update B, C, D, E ---

$ws.Range("B32").Value = 9000041
$ws.Range("C32").Value = 8000041
$ws.Range("D32").Value = "IRPI"
$ws.Range("E32").Value = "Modelli e carte di suscettibilità da frana"

$ws.Range("B33").Value = 9000045
$ws.Range("C33").Value = 8000042
$ws.Range("D33").Value = "Confederazion elvetica"
$ws.Range("E33").Value = "Legge federale sulla protezione della popolazione e sulla protezione civile del 4 ottobre 2002"

$ws.Range("B34").Value = 9000052
$ws.Range("C34").Value = 8000047
$ws.Range("D34").Value = "Ufficio Federale della Protezione della Popolazione"
$ws.Range("E34").Value = "Legge sulla protezione della popolazione (del 26 febbraio 2007)"

$ws.Range("B35").Value = 9000052
$ws.Range("C35").Value = 8000052
$ws.Range("E35").Value = "Il comando della protezione civile - Personale"

$ws.Range("B36").Value = 9000062
$ws.Range("C36").Value = 8000062
$ws.Range("D36").Value = "Repubblica e Cantone Ticino"
$ws.Range("E36").Value = "Servizio della protezione della popolazione"

$ws.Range("B37").Value = 9000065
$ws.Range("C37").Value = 8000065
$ws.Range("D37").Value = "Tommaso Sansone"
$ws.Range("E37").Value = "-- documento confronto normativa -- wp 3.2 gestisco -- da completare"

$ws.Range("B38").Value = 9000065
$ws.Range("C38").Value = 8000070
$ws.Range("E38").Value = "-- -- documento confronto normativa -- wp 3.2 gestisco -- da completare"

$ws.Range("B39").Value = 9000067
$ws.Range("C39").Value = 8000065
$ws.Range("D39").Value = "Tommaso sansone"
$ws.Range("E39").Value = "-- documento confronto normativa -- wp 3.2 gestisco -- da completare"

$ws.Range("C40").Value = 8000032
$ws.Range("D40").Value = "Alberto Bruno, Funzionario della protezione civile di regione lombardia"
$ws.Range("E40").Value = ""
